# Generate Report for Handback
# Update the timestamp strings recorded for the handback/handoff report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview.Range("G2").Value = "2016-08-15 17:01:24"

# zh-cn sheet: "Correspond Handoff Datetime" (column H) / "Correspond Handback DateTime" (column K)
$wsZhCn.Range("H2").Value = "2016-08-15 17:01:18"
$wsZhCn.Range("K2").Value = "2016-08-15 17:01:36"

# de-de sheet: "Correspond Handoff Datetime" (column H) / "Correspond Handback DateTime" (column K)
$wsDeDe.Range("H2").Value = "2016-08-15 17:01:24"
$wsDeDe.Range("K2").Value = "2016-08-15 17:01:43"
